# Correct unit of investments from mio€ to € by multiplying the data
# range (C2:M18) by 1,000,000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:M18")
foreach ($cell in $rng.Cells) {
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current * 1000000
    }
}

$ws.Range("C19").Select()
